$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Rows 117-120: data was rotated (row 117<-old119, 118<-old120,
# 119<-old118, 120<-old117) while the index column A stays fixed.
# Apply the resulting values explicitly, column by column.
# -----------------------------------------------------------------

# Row 117 (was old row 119's data)
$ws.Range("B117").Value = 7013702
$ws.Range("C117").Value = "Uruguay Primera División"
$ws.Range("D117").Value = 45267.70833333334
$ws.Range("E117").Value = "Defensor Sporting"
$ws.Range("F117").Value = "Danubio"
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 2
$ws.Range("I117").Value = "A"
$ws.Range("J117").Value = 1.8
$ws.Range("K117").Value = 3.6
$ws.Range("L117").Value = 4.2
$ws.Range("M117").Value = 1.8
$ws.Range("N117").Value = 3.6
$ws.Range("O117").Value = 4.2
$ws.Range("P117").Value = -0.75
$ws.Range("Q117").Value = 2.05
$ws.Range("R117").Value = 1.8
$ws.Range("S117").Value = 2.25
$ws.Range("T117").Value = 1.85
$ws.Range("U117").Value = 2
$ws.Range("V117").Value = -1
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 3.2
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = 0.8
$ws.Range("AA117").Value = -0.5
$ws.Range("AB117").Value = 0.5

# Row 118 (was old row 120's data)
$ws.Range("B118").Value = 7013885
$ws.Range("C118").Value = "Uruguay Primera División"
$ws.Range("D118").Value = 45267.70833333334
$ws.Range("E118").Value = "La Luz"
$ws.Range("F118").Value = "Atletico Fenix Montevideo"
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 2
$ws.Range("I118").Value = "A"
$ws.Range("J118").Value = 3
$ws.Range("K118").Value = 3
$ws.Range("L118").Value = 2.4
$ws.Range("M118").Value = 2.9
$ws.Range("N118").Value = 2.75
$ws.Range("O118").Value = 2.6
$ws.Range("P118").Value = 0
$ws.Range("Q118").Value = 2.025
$ws.Range("R118").Value = 1.825
$ws.Range("S118").Value = 2
$ws.Range("T118").Value = 2.025
$ws.Range("U118").Value = 1.825
$ws.Range("V118").Value = -1
$ws.Range("W118").Value = -1
$ws.Range("X118").Value = 1.6
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 0.825
$ws.Range("AA118").Value = 0
$ws.Range("AB118").Value = 0

# Row 119 (was old row 118's data)
$ws.Range("B119").Value = 7013409
$ws.Range("C119").Value = "Uruguay Primera División"
$ws.Range("D119").Value = 45267.70833333334
$ws.Range("E119").Value = "Nacional De Football"
$ws.Range("F119").Value = "Torque"
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 1
$ws.Range("I119").Value = "D"
$ws.Range("J119").Value = 1.666
$ws.Range("K119").Value = 3.9
$ws.Range("L119").Value = 4.5
$ws.Range("M119").Value = 1.615
$ws.Range("N119").Value = 4
$ws.Range("O119").Value = 4.75
$ws.Range("P119").Value = -0.75
$ws.Range("Q119").Value = 1.8
$ws.Range("R119").Value = 2.05
$ws.Range("S119").Value = 2.75
$ws.Range("T119").Value = 1.95
$ws.Range("U119").Value = 1.9
$ws.Range("V119").Value = -1
$ws.Range("W119").Value = 3
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = 1.05
$ws.Range("AA119").Value = -1
$ws.Range("AB119").Value = 0.8999999999999999

# Row 120 (was old row 117's data)
$ws.Range("B120").Value = 7013886
$ws.Range("C120").Value = "Uruguay Primera División"
$ws.Range("D120").Value = 45267.70833333334
$ws.Range("E120").Value = "Racing Club de Montevideo"
$ws.Range("F120").Value = "Cerro"
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1
$ws.Range("I120").Value = "A"
$ws.Range("J120").Value = 2.25
$ws.Range("K120").Value = 3.1
$ws.Range("L120").Value = 3.25
$ws.Range("M120").Value = 2.25
$ws.Range("N120").Value = 2.875
$ws.Range("O120").Value = 3.5
$ws.Range("P120").Value = -0.25
$ws.Range("Q120").Value = 1.95
$ws.Range("R120").Value = 1.9
$ws.Range("S120").Value = 2
$ws.Range("T120").Value = 1.925
$ws.Range("U120").Value = 1.925
$ws.Range("V120").Value = -1
$ws.Range("W120").Value = -1
$ws.Range("X120").Value = 2.5
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 0.8999999999999999
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = 0.925

# -----------------------------------------------------------------
# Odds re-calculation tweaks on rows 226, 227, 229, 231, 232
# -----------------------------------------------------------------

# Row 226
$ws.Range("M226").Value = 1.5
$ws.Range("O226").Value = 6.5
$ws.Range("Q226").Value = 1.85
$ws.Range("R226").Value = 2
$ws.Range("T226").Value = 1.875
$ws.Range("U226").Value = 1.975

# Row 227
$ws.Range("Q227").Value = 1.95
$ws.Range("R227").Value = 1.9

# Row 229
$ws.Range("Q229").Value = 1.85
$ws.Range("R229").Value = 2

# Row 231
$ws.Range("Q231").Value = 2.025
$ws.Range("R231").Value = 1.825

# Row 232
$ws.Range("T232").Value = 1.9
$ws.Range("U232").Value = 1.95
